{"js": "// Three small wording tweaks in the script dialogue:\n//   1) \"...but maybe his parent is inside?\" -> \"...but maybe his parent's inside?\"\n//   2) \"Well, even if his parent is inside the store...\" -> \"Well, even if his parent's inside the store...\"\n//   3) \"?Greta (embarrassed embarrassed): I\u2026 uh\u2026\" -> \"?Greta (neutral embarrassed): I\u2026 uh\u2026\"\n\nconst RIGHT_SINGLE_QUOTE = \"\\u2019\";\n\n// 1) \"his parent is inside?\" -> \"his parent's inside?\"\nconst search1 = context.document.body.search(\"his parent is inside?\", { matchCase: true });\nsearch1.load(\"text\");\nawait context.sync();\nif (search1.items.length > 0) {\n  search1.items[0].insertText(\"his parent\" + RIGHT_SINGLE_QUOTE + \"s inside?\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) \"his parent is inside the store\" -> \"his parent's inside the store\"\nconst search2 = context.document.body.search(\"his parent is inside the store\", { matchCase: true });\nsearch2.load(\"text\");\nawait context.sync();\nif (search2.items.length > 0) {\n  search2.items[0].insertText(\"his parent\" + RIGHT_SINGLE_QUOTE + \"s inside the store\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 3) \"(embarrassed embarrassed)\" -> \"(neutral embarrassed)\"\nconst search3 = context.document.body.search(\"(embarrassed embarrassed)\", { matchCase: true });\nsearch3.load(\"text\");\nawait context.sync();\nif (search3.items.length > 0) {\n  search3.items[0].insertText(\"(neutral embarrassed)\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Three small wording tweaks in the script dialogue:\n#   1) \"...but maybe his parent is inside?\" -> \"...but maybe his parent's inside?\"\n#   2) \"Well, even if his parent is inside the store...\" -> \"Well, even if his parent's inside the store...\"\n#   3) \"?Greta (embarrassed embarrassed): I... uh...\" -> \"?Greta (neutral embarrassed): I... uh...\"\n\n$d = $word.ActiveDocument\n$RSQUO = [char]0x2019\n\n# wdReplace constants\n$wdFindContinue = 1\n$wdReplaceOne = 1\n\n# 1) \"his parent is inside?\" -> \"his parent's inside?\"\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"his parent is inside?\"\n$find.Replacement.Text = \"his parent\" + $RSQUO + \"s inside?\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceOne) | Out-Null\n\n# 2) \"his parent is inside the store\" -> \"his parent's inside the store\"\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"his parent is inside the store\"\n$find.Replacement.Text = \"his parent\" + $RSQUO + \"s inside the store\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceOne) | Out-Null\n\n# 3) \"(embarrassed embarrassed)\" -> \"(neutral embarrassed)\"\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"(embarrassed embarrassed)\"\n$find.Replacement.Text = \"(neutral embarrassed)\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceOne) | Out-Null\n"}
